$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 151, shifting existing rows 151:182 down to 152:183
$ws.Rows.Item(151).Insert()

# Populate the new row 151 with the latest weekly price record
$ws.Cells.Item(151, 1).Value = 7
$ws.Cells.Item(151, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(151, 3).Value = "Ñuble"
$ws.Cells.Item(151, 4).Value = 44476
$ws.Cells.Item(151, 5).Value = 16
$ws.Cells.Item(151, 6).Value = 100114013
$ws.Cells.Item(151, 7).Value = "Zanahoria"
$ws.Cells.Item(151, 8).Value = "Sin especificar"
$ws.Cells.Item(151, 9).Value = "Primera"
$ws.Cells.Item(151, 10).Value = 120
$ws.Cells.Item(151, 11).Value = 8000
$ws.Cells.Item(151, 12).Value = 9000
$ws.Cells.Item(151, 13).Value = 8500
$ws.Cells.Item(151, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(151, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(151, 16).Value = 425
$ws.Cells.Item(151, 17).Value = 20
$ws.Cells.Item(151, 18).Value = "Hortaliza"
